# Apply updated dSF (column F) values to the data rows on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = 3
    5  = 0
    7  = 6
    9  = 3
    10 = -3
    11 = -2
    12 = -4
    13 = -2
    15 = -1
    16 = -4
    17 = -1
    18 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
